# Daily attendance processing - 2025-11-15 06:29:17
#
# For every data row, the "Recorded By" column (G) holds a comma
# separated list of recorder names/emails (e.g. "dnasr281@gmail.com,
# System"). Re-sort each list into (ordinal/ASCII) alphabetical order,
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".

# Pads a numeric "char code" key string out to a fixed length so that
# plain lexicographic comparisons of two keys behave like comparing the
# original strings position by position (shorter strings sort first).
function PadKey($k) {
    $target = 200
    while ($k.Length -lt $target) {
        $k = [string]::Concat($k, "0")
    }
    return $k
}

# Builds an ordinal sort key for a string: each character becomes a
# zero padded 5-digit code, concatenated together, then padded to a
# fixed width. Comparing these keys as plain strings reproduces a
# case-sensitive, ordinal (ASCII) string comparison.
function MakeKey($s) {
    $key = ""
    $len = $s.Length
    for ($i = 0; $i -lt $len; $i++) {
        $code = [int][char]$s[$i]
        $piece = [string]::Format("{0:D5}", $code)
        $key = [string]::Concat($key, $piece)
    }
    return PadKey $key
}

# Splits a comma separated cell value, trims each entry, sorts the
# entries in ordinal order and re-joins them with ", ".
function SortCsv($value) {
    $parts = $value -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $objs = @()
    foreach ($t in $trimmed) {
        $objs += @{ Text = $t; Key = (MakeKey $t) }
    }

    $sorted = $objs | Sort-Object -Property Key

    $result = @()
    foreach ($o in $sorted) {
        $result += $o.Text
    }

    return [string]::Join(", ", $result)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G is the "Recorded By" column.
$col = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $value = $cell.Value()

    if ($value -ne $null -and $value -ne "") {
        $newValue = SortCsv $value
        if ($newValue -ne $value) {
            $cell.Value = $newValue
        }
    }
}
